$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.814.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2560"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07698"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.649.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.856.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5419"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7880"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.830.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "200.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.320"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.867"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.931"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.929"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1132"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.682"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.238"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04978"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.264"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.534"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.368"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.164.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.622"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5554"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.656"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8015"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.768.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4513"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9982"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05066"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "
